# Add new joinee rows to the database sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of IDs (col A) and values (col B) appended after row 555.
$newRows = @(
    @("21F1000164", 291),
    @("21F1000925", 264),
    @("21F1001180", 264),
    @("21F1001914", 311),
    @("21F1002604", 165),
    @("21F1002644", 264),
    @("21F1003008", 311),
    @("21F1003057", 311),
    @("21F1003406", 311),
    @("21F1003511", 264),
    @("21F1004030", 264),
    @("21F1004033", 311),
    @("21F1004833", 311),
    @("21F1005277", 311),
    @("21F1005925", 311),
    @("21F2000104", 311)
)

$startRow = 556
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $id = $newRows[$i][0]
    $val = $newRows[$i][1]
    $ws.Cells.Item($r, 1).Value = $id
    $ws.Cells.Item($r, 2).Value = $val
}

$lastRow = $startRow + $newRows.Length - 1

# Update the visible selection to span the newly extended data range.
$ws.Range("A4:A$lastRow").Select() | Out-Null

# Zoom level was reduced from 160% to 140%.
$ws.Application.ActiveWindow.Zoom = 140

# Column widths were tweaked slightly.
$ws.Columns.Item(1).ColumnWidth = 41.0
$ws.Columns.Item(2).ColumnWidth = 6.165
